$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 14

$ws.Cells.Item($row, 1).Value = 42619.89266203704
$ws.Cells.Item($row, 2).Value = 26
$ws.Cells.Item($row, 3).Value = 65
$ws.Cells.Item($row, 4).Value = 34
$ws.Cells.Item($row, 5).Value = 65
$ws.Cells.Item($row, 6).Value = 35
$ws.Cells.Item($row, 7).Value = 23315
$ws.Cells.Item($row, 8).Value = 19358
$ws.Cells.Item($row, 9).Value = 1105
$ws.Cells.Item($row, 10).Value = 251
$ws.Cells.Item($row, 11).Value = 131
$ws.Cells.Item($row, 12).Value = 20
$ws.Cells.Item($row, 13).Value = 11
$ws.Cells.Item($row, 14).Value = "Named"

$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"
